$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3322.4138
$ws.Range("I64").Value = 2927.1428
$ws.Range("K64").Value = 2927.1428
$ws.Range("M64").Value = -2679.1428
$ws.Range("H67").Value = 3322.4138
$ws.Range("I67").Value = 2927.1428
$ws.Range("K67").Value = 2927.1428
$ws.Range("M67").Value = -2069.1428
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H98").Value = 1722.8572
$ws.Range("I98").Value = 1701.1538
$ws.Range("J98").Value = 2005
$ws.Range("K98").Value = 1701.1538
$ws.Range("L98").Value = 2005
$ws.Range("M98").Value = -203.1538
$ws.Range("N98").Value = -5001
$ws.Range("H116").Value = 3437.25
$ws.Range("I116").Value = 2749.1667
$ws.Range("K116").Value = 2749.1667
$ws.Range("M116").Value = 692.8332999999998
$ws.Range("H122").Value = 1722.8572
$ws.Range("I122").Value = 1701.1538
$ws.Range("J122").Value = 2005
$ws.Range("K122").Value = 5103.4614
$ws.Range("L122").Value = 6015
$ws.Range("M122").Value = -2653.4614
$ws.Range("N122").Value = -10915
$ws.Range("H137").Value = 1107.8158
$ws.Range("I137").Value = 960.8
$ws.Range("J137").Value = 1390.5385
$ws.Range("K137").Value = 2882.4
$ws.Range("L137").Value = 4171.6155
$ws.Range("M137").Value = -332.3999999999996
$ws.Range("N137").Value = -9271.6155
$ws.Range("H138").Value = 2137.57
$ws.Range("I138").Value = 932.0909
$ws.Range("J138").Value = 2286.5618
$ws.Range("K138").Value = 2796.2727
$ws.Range("L138").Value = 6859.6854
$ws.Range("M138").Value = 2343.7273
$ws.Range("N138").Value = -17139.6854

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 2999
$ws.Range("I25").Value = 2999
$ws.Range("K25").Value = 2999
$ws.Range("M25").Value = -2597
$ws.Range("H32").Value = 4351.35
$ws.Range("I32").Value = 4830.514
$ws.Range("J32").Value = 997.2
$ws.Range("K32").Value = 4830.514
$ws.Range("L32").Value = 997.2
$ws.Range("M32").Value = -4543.514
$ws.Range("N32").Value = -1571.2
$ws.Range("H45").Value = 2008
$ws.Range("I45").Value = 1954.2222
$ws.Range("J45").Value = 2250
$ws.Range("K45").Value = 1954.2222
$ws.Range("L45").Value = 2250
$ws.Range("M45").Value = -1577.2222
$ws.Range("N45").Value = -3004
$ws.Range("H61").Value = 1323.8
$ws.Range("I61").Value = 1110.125
$ws.Range("J61").Value = 2178.5
$ws.Range("K61").Value = 1110.125
$ws.Range("L61").Value = 2178.5
$ws.Range("M61").Value = -898.125
$ws.Range("N61").Value = -2602.5
$ws.Range("H74").Value = 1097.5
$ws.Range("I74").Value = 1102.4615
$ws.Range("K74").Value = 1102.4615
$ws.Range("M74").Value = -228.4614999999999
$ws.Range("H77").Value = 1097.5
$ws.Range("I77").Value = 1102.4615
$ws.Range("K77").Value = 5512.307499999999
$ws.Range("M77").Value = -1144.307499999999
$ws.Range("H97").Value = 482
$ws.Range("I97").Value = 482
$ws.Range("K97").Value = 482
$ws.Range("M97").Value = 14
$ws.Range("H132").Value = 2998.9312
$ws.Range("I132").Value = 2821.9092
$ws.Range("K132").Value = 8465.7276
$ws.Range("M132").Value = -5935.7276
$ws.Range("H136").Value = 1323.8
$ws.Range("I136").Value = 1110.125
$ws.Range("J136").Value = 2178.5
$ws.Range("K136").Value = 3330.375
$ws.Range("L136").Value = 6535.5
$ws.Range("M136").Value = -780.375
$ws.Range("N136").Value = -11635.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 50001020
$ws.Range("J94").Value = 3000
$ws.Range("L94").Value = 3000
$ws.Range("N94").Value = -3902
$ws.Range("H105").Value = 55557644
$ws.Range("I105").Value = 83335120
$ws.Range("J105").Value = 2685.1667
$ws.Range("K105").Value = 83335120
$ws.Range("L105").Value = 2685.1667
$ws.Range("M105").Value = -83333373
$ws.Range("N105").Value = -6179.1667
$ws.Range("H134").Value = 5891.3477
$ws.Range("I134").Value = 1568.7059
$ws.Range("J134").Value = 18138.834
$ws.Range("K134").Value = 4706.1177
$ws.Range("L134").Value = 54416.50199999999
$ws.Range("M134").Value = -2171.1177
$ws.Range("N134").Value = -59486.50199999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 780.95654
$ws.Range("I31").Value = 715.8889
$ws.Range("J31").Value = 902.9583
$ws.Range("K31").Value = 715.8889
$ws.Range("L31").Value = 902.9583
$ws.Range("M31").Value = -420.8889
$ws.Range("N31").Value = -1492.9583
$ws.Range("H34").Value = 780.95654
$ws.Range("I34").Value = 715.8889
$ws.Range("J34").Value = 902.9583
$ws.Range("K34").Value = 715.8889
$ws.Range("L34").Value = 902.9583
$ws.Range("M34").Value = -513.8889
$ws.Range("N34").Value = -1306.9583
$ws.Range("H99").Value = 1790.2941
$ws.Range("I99").Value = 1748.7273
$ws.Range("K99").Value = 1748.7273
$ws.Range("M99").Value = -250.7273
$ws.Range("H126").Value = 1790.2941
$ws.Range("I126").Value = 1748.7273
$ws.Range("K126").Value = 5246.1819
$ws.Range("M126").Value = -2776.1819
$ws.Range("H132").Value = 12501.1
$ws.Range("I132").Value = 18668.666
$ws.Range("J132").Value = 3249.75
$ws.Range("K132").Value = 56005.99800000001
$ws.Range("L132").Value = 9749.25
$ws.Range("M132").Value = -53475.99800000001
$ws.Range("N132").Value = -14809.25
$ws.Range("H134").Value = 20835514
$ws.Range("I134").Value = 30304984
$ws.Range("J134").Value = 2680
$ws.Range("K134").Value = 90914952
$ws.Range("L134").Value = 8040
$ws.Range("M134").Value = -90912417
$ws.Range("N134").Value = -13110

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 45456480
$ws.Range("J131").Value = 2500.4
$ws.Range("L131").Value = 7501.200000000001
$ws.Range("N131").Value = -17581.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 19016.5
$ws.Range("J95").Value = 19016.5
$ws.Range("L95").Value = 19016.5
$ws.Range("N95").Value = -24508.5
$ws.Range("H132").Value = 2205.3447
$ws.Range("I132").Value = 1776.5
$ws.Range("K132").Value = 5329.5
$ws.Range("M132").Value = -2799.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1891.5
$ws.Range("I7").Value = 1925
$ws.Range("J7").Value = 1874.75
$ws.Range("K7").Value = 1925
$ws.Range("L7").Value = 1874.75
$ws.Range("M7").Value = -1813
$ws.Range("N7").Value = -2098.75
$ws.Range("H22").Value = 1287.375
$ws.Range("J22").Value = 1699.8
$ws.Range("L22").Value = 1699.8
$ws.Range("N22").Value = -2289.8
$ws.Range("H27").Value = 1287.375
$ws.Range("J27").Value = 1699.8
$ws.Range("L27").Value = 1699.8
$ws.Range("N27").Value = -1913.8
$ws.Range("H40").Value = 2699.3
$ws.Range("I40").Value = 2284.8572
$ws.Range("K40").Value = 2284.8572
$ws.Range("M40").Value = -2148.8572
$ws.Range("H106").Value = 28941
$ws.Range("J106").Value = 28941
$ws.Range("L106").Value = 28941
$ws.Range("N106").Value = -31465
$ws.Range("H112").Value = 47570.43
$ws.Range("J112").Value = 47570.43
$ws.Range("L112").Value = 47570.43
$ws.Range("N112").Value = -50524.43
$ws.Range("H126").Value = 1891.5
$ws.Range("I126").Value = 1925
$ws.Range("J126").Value = 1874.75
$ws.Range("K126").Value = 5775
$ws.Range("L126").Value = 5624.25
$ws.Range("M126").Value = -3305
$ws.Range("N126").Value = -10564.25
$ws.Range("H132").Value = 40721.46
$ws.Range("I132").Value = 1953.2307
$ws.Range("J132").Value = 79489.69500000001
$ws.Range("K132").Value = 5859.6921
$ws.Range("L132").Value = 238469.085
$ws.Range("M132").Value = -3329.6921
$ws.Range("N132").Value = -243529.085
$ws.Range("H136").Value = 2442.7144
$ws.Range("I136").Value = 1924.75
$ws.Range("J136").Value = 3133.3333
$ws.Range("K136").Value = 5774.25
$ws.Range("L136").Value = 9399.999899999999
$ws.Range("M136").Value = -3224.25
$ws.Range("N136").Value = -14499.9999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 76924460
$ws.Range("J126").Value = 2245.6
$ws.Range("L126").Value = 6736.799999999999
$ws.Range("N126").Value = -11676.8
$ws.Range("H132").Value = 2373.6365
$ws.Range("I132").Value = 2201.4
$ws.Range("K132").Value = 6604.200000000001
$ws.Range("M132").Value = -4074.200000000001
